$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update E6 value (0 -> 0.5)
$ws.Range("E6").Value = 0.5

# Add new label in G6 referencing new shared string "min error formula"
$ws.Range("G6").Value = "min error formula"

# Update the selected cell from E7 to G7
$ws.Range("G7").Select()
